$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 997.53656
$ws.Range("I18").Value = 959.975
$ws.Range("J18").Value = 2500
$ws.Range("K18").Value = 959.975
$ws.Range("L18").Value = 2500
$ws.Range("M18").Value = -675.975
$ws.Range("N18").Value = -3068
$ws.Range("H40").Value = 4862.375
$ws.Range("I40").Value = 9666.333000000001
$ws.Range("J40").Value = 1980
$ws.Range("K40").Value = 9666.333000000001
$ws.Range("L40").Value = 1980
$ws.Range("M40").Value = -9491.333000000001
$ws.Range("N40").Value = -2330
$ws.Range("H41").Value = 1296.2
$ws.Range("I41").Value = 1399.6666
$ws.Range("K41").Value = 1399.6666
$ws.Range("M41").Value = -959.6666
$ws.Range("H95").Value = 48820.89
$ws.Range("J95").Value = 48054.75
$ws.Range("L95").Value = 48054.75
$ws.Range("N95").Value = -53546.75
$ws.Range("H98").Value = 3320.25
$ws.Range("I98").Value = 2810.394
$ws.Range("J98").Value = 4849.8184
$ws.Range("K98").Value = 2810.394
$ws.Range("L98").Value = 4849.8184
$ws.Range("M98").Value = -1312.394
$ws.Range("N98").Value = -7845.8184
$ws.Range("H106").Value = 6377.1055
$ws.Range("I106").Value = 4342
$ws.Range("J106").Value = 12075.4
$ws.Range("K106").Value = 4342
$ws.Range("L106").Value = 12075.4
$ws.Range("M106").Value = -3711
$ws.Range("N106").Value = -13337.4
$ws.Range("H122").Value = 3320.25
$ws.Range("I122").Value = 2810.394
$ws.Range("J122").Value = 4849.8184
$ws.Range("K122").Value = 8431.181999999999
$ws.Range("L122").Value = 14549.4552
$ws.Range("M122").Value = -5981.181999999999
$ws.Range("N122").Value = -19449.4552
$ws.Range("H132").Value = 9078.315000000001
$ws.Range("I132").Value = 6837.029
$ws.Range("J132").Value = 35226.668
$ws.Range("K132").Value = 20511.087
$ws.Range("L132").Value = 105680.004
$ws.Range("M132").Value = -17981.087
$ws.Range("N132").Value = -110740.004
$ws.Range("H135").Value = 2246.4
$ws.Range("I135").Value = 2016.125
$ws.Range("J135").Value = 3167.5
$ws.Range("K135").Value = 18145.125
$ws.Range("L135").Value = 28507.5
$ws.Range("M135").Value = -15610.125
$ws.Range("N135").Value = -33577.5
$ws.Range("H138").Value = 2199.76
$ws.Range("I138").Value = 1466.069
$ws.Range("J138").Value = 2499.4365
$ws.Range("K138").Value = 4398.207
$ws.Range("L138").Value = 7498.309499999999
$ws.Range("M138").Value = 741.7929999999997
$ws.Range("N138").Value = -17778.3095

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 6413.6665
$ws.Range("I28").Value = 6413.6665
$ws.Range("K28").Value = 6413.6665
$ws.Range("M28").Value = -6221.6665
$ws.Range("H32").Value = 4739.662
$ws.Range("J32").Value = 14791.23
$ws.Range("L32").Value = 14791.23
$ws.Range("N32").Value = -15365.23
$ws.Range("H94").Value = 46833.332
$ws.Range("J94").Value = 46833.332
$ws.Range("L94").Value = 46833.332
$ws.Range("N94").Value = -48635.332
$ws.Range("H99").Value = 6413.6665
$ws.Range("I99").Value = 6413.6665
$ws.Range("K99").Value = 6413.6665
$ws.Range("M99").Value = -3418.6665
$ws.Range("H108").Value = 53623.09
$ws.Range("J108").Value = 53623.09
$ws.Range("L108").Value = 53623.09
$ws.Range("N108").Value = -61303.09

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 25382.182
$ws.Range("J99").Value = 23275.8
$ws.Range("L99").Value = 23275.8
$ws.Range("N99").Value = -26271.8
$ws.Range("H107").Value = 1536.1923
$ws.Range("I107").Value = 1589.25
$ws.Range("K107").Value = 1589.25
$ws.Range("M107").Value = 330.75
$ws.Range("H134").Value = 57817.652
$ws.Range("I134").Value = 84012.766
$ws.Range("J134").Value = 23764
$ws.Range("K134").Value = 252038.298
$ws.Range("L134").Value = 71292
$ws.Range("M134").Value = -249503.298
$ws.Range("N134").Value = -76362

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 864.5599999999999
$ws.Range("I22").Value = 428.83334
$ws.Range("K22").Value = 428.83334
$ws.Range("M22").Value = -78.83334000000002
$ws.Range("H31").Value = 14828.219
$ws.Range("I31").Value = 6183.12
$ws.Range("J31").Value = 45703.57
$ws.Range("K31").Value = 6183.12
$ws.Range("L31").Value = 45703.57
$ws.Range("M31").Value = -5888.12
$ws.Range("N31").Value = -46293.57
$ws.Range("H34").Value = 14828.219
$ws.Range("I34").Value = 6183.12
$ws.Range("J34").Value = 45703.57
$ws.Range("K34").Value = 6183.12
$ws.Range("L34").Value = 45703.57
$ws.Range("M34").Value = -5981.12
$ws.Range("N34").Value = -46107.57
$ws.Range("H41").Value = 4773.6
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("H50").Value = 30000
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("H59").Value = 24784.666
$ws.Range("J59").Value = 44250
$ws.Range("L59").Value = 44250
$ws.Range("N59").Value = -46540
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("H132").Value = 4473.057
$ws.Range("I132").Value = 1375.5807
$ws.Range("J132").Value = 28478.5
$ws.Range("K132").Value = 4126.742099999999
$ws.Range("L132").Value = 85435.5
$ws.Range("M132").Value = -1596.742099999999
$ws.Range("N132").Value = -90495.5
$ws.Range("H134").Value = 3599.742
$ws.Range("I134").Value = 1486.4468
$ws.Range("K134").Value = 4459.3404
$ws.Range("M134").Value = -1924.3404
$ws.Range("N41").ClearContents()
$ws.Range("N50").ClearContents()
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 503
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("H80").Value = 13423.5
$ws.Range("I80").Value = 10231.667
$ws.Range("K80").Value = 30695.001
$ws.Range("M80").Value = -29759.001
$ws.Range("H83").Value = 13423.5
$ws.Range("I83").Value = 10231.667
$ws.Range("K83").Value = 92085.003
$ws.Range("M83").Value = -87405.003
$ws.Range("H86").Value = 790.2857
$ws.Range("J86").Value = 899
$ws.Range("L86").Value = 2697
$ws.Range("N86").Value = -5069
$ws.Range("H89").Value = 790.2857
$ws.Range("J89").Value = 899
$ws.Range("L89").Value = 8091
$ws.Range("N89").Value = -19947
$ws.Range("H122").Value = 9360138
$ws.Range("J122").Value = 1671097.9
$ws.Range("L122").Value = 15039881.1
$ws.Range("N122").Value = -15044781.1
$ws.Range("H129").Value = 2071.1428
$ws.Range("I129").Value = 1043
$ws.Range("K129").Value = 3129
$ws.Range("M129").Value = 1871
$ws.Range("N47").ClearContents()

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 12500
$ws.Range("J18").Value = 12500
$ws.Range("L18").Value = 12500
$ws.Range("N18").Value = -13086
$ws.Range("H44").Value = 16888.445
$ws.Range("I44").Value = 20499.834
$ws.Range("J44").Value = 9665.666999999999
$ws.Range("K44").Value = 20499.834
$ws.Range("L44").Value = 9665.666999999999
$ws.Range("M44").Value = -19903.834
$ws.Range("N44").Value = -10857.667
$ws.Range("H102").Value = 1019442
$ws.Range("I102").Value = 1448631.6
$ws.Range("K102").Value = 1448631.6
$ws.Range("M102").Value = -1447009.6
$ws.Range("H132").Value = 9535.772000000001
$ws.Range("I132").Value = 7837.6924
$ws.Range("K132").Value = 23513.0772
$ws.Range("M132").Value = -20983.0772

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 90912140
$ws.Range("I16").Value = 90912140
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 90912140
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -90911970
$ws.Range("H46").Value = 4199.6665
$ws.Range("I46").Value = 4199.6665
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 4199.6665
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -4011.6665
$ws.Range("H93").Value = 28578304
$ws.Range("I93").Value = 45459436
$ws.Range("J93").Value = 10233.77
$ws.Range("K93").Value = 45459436
$ws.Range("L93").Value = 10233.77
$ws.Range("M93").Value = -45458188
$ws.Range("N93").Value = -12729.77
$ws.Range("H100").Value = 5837.56
$ws.Range("I100").Value = 4219.722
$ws.Range("J100").Value = 9997.714
$ws.Range("K100").Value = 4219.722
$ws.Range("L100").Value = 9997.714
$ws.Range("M100").Value = -3678.722
$ws.Range("N100").Value = -11079.714
$ws.Range("H103").Value = 17520.4
$ws.Range("J103").Value = 17520.4
$ws.Range("L103").Value = 17520.4
$ws.Range("N103").Value = -19864.4
$ws.Range("N16").ClearContents()
$ws.Range("N46").ClearContents()

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 41000
$ws.Range("J75").Value = 41000
$ws.Range("L75").Value = 41000
$ws.Range("N75").Value = -42872
$ws.Range("H78").Value = 41000
$ws.Range("J78").Value = 41000
$ws.Range("L78").Value = 123000
$ws.Range("N78").Value = -132360
$ws.Range("H97").Value = 46576.668
$ws.Range("J97").Value = 46576.668
$ws.Range("L97").Value = 46576.668
$ws.Range("N97").Value = -48558.668
